$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 93 (shifts old rows 93..113 down to 94..114) ---
$ws.Rows.Item(93).Insert()

# Copy formatting (style) of the date cell above into the newly inserted date cell
# so A93 keeps the same date number-format / border as the rest of column A.
$ws.Cells.Item(92,1).Copy()
$ws.Cells.Item(93,1).PasteSpecial(-4122)

# --- 2. Fill in the values for the new row 93 ---
$ws.Cells.Item(93,1).Value = 44235
$ws.Cells.Item(93,2).Value = 1
$ws.Cells.Item(93,3).Value = 10
$ws.Cells.Item(93,4).Value = 160.2307322544464

# --- 3. Rolling-window values recomputed for rows 90, 91 (swap) and 94-96 ---
$ws.Cells.Item(90,3).Value = 12
$ws.Cells.Item(90,4).Value = 192.2768787053357

$ws.Cells.Item(91,3).Value = 13
$ws.Cells.Item(91,4).Value = 208.2999519307803

$ws.Cells.Item(94,3).Value = 10
$ws.Cells.Item(94,4).Value = 160.2307322544464

$ws.Cells.Item(95,3).Value = 10
$ws.Cells.Item(95,4).Value = 160.2307322544464

$ws.Cells.Item(96,3).Value = 9
$ws.Cells.Item(96,4).Value = 144.2076590290018

# --- 4. Row 112 (old row 111, shifted) now gets its rolling-window values filled in ---
$ws.Cells.Item(112,3).Value = 10
$ws.Cells.Item(112,4).Value = 160.2307322544464

# --- 5. Append two brand-new rows (114, 115) at the end of the series ---
$ws.Cells.Item(113,1).Copy()
$ws.Range("A114:A115").PasteSpecial(-4122)

$ws.Cells.Item(114,1).Value = 44256
$ws.Cells.Item(114,2).Value = 0

$ws.Cells.Item(115,1).Value = 44257
$ws.Cells.Item(115,2).Value = 4

Write-Host ("UsedRange=" + $ws.UsedRange.Address())
